# Update "New Microsoft Excel Worksheet.xlsx":
#  - enter 11 / 111 into I11 / I12 on Sheet1
#  - leave the active-cell selection on I13 (the cell below the entered data)
#  - switch the workbook's default/Normal font from Calibri to Arial

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the two numeric values that become the sheet's new used range (I11:I12)
$ws.Range("I11").Value = 11
$ws.Range("I12").Value = 111

# After typing into I12 and pressing Enter, Excel leaves the active cell on I13 -
# reproduce that saved selection state.
$ws.Range("I13").Select()

# Change the workbook's default ("Normal") font from Calibri to Arial - this is
# the font used by every cell that has no explicit font override, i.e. I11/I12
# above keep using the (now Arial) default style rather than an explicit font.
$wb.Styles("Normal").Font.Name = "Arial"
